$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portugal Liga 3")
$ws.Cells.Item(3,11).Value2 = "A"
